# petty-cashBook-2021.xlsx — 12-Apr-2021 midday update
# The daily cash book (Sheet1) is rolled forward: the old entries for
# 9-12 Feb 2021 (rows 3-32) are cleared out and the running balance
# (SALDO AWAL, E2) is reset to the new carried-forward amount. The first
# row of the new period (row 3) is re-dated and its Debit entry is blank
# until postings resume.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New opening balance (SALDO AWAL) carried forward
$ws.Range("E2").Value = 152525

# Re-date the first transaction row to the new period
$ws.Range("A3").Value = 44298

# Clear the old Debit figure on row 3 (no postings yet this period)
$ws.Range("D3").Clear()

# Wipe out all of the old transaction detail (dates/descriptions/amounts)
# for the rest of the rolled-over block; the running-balance formulas in
# column E stay in place and simply recompute to the flat new balance.
$ws.Range("A4:D32").Clear()

# Leave the cursor where the midday edit left off
[void]$ws.Activate()
[void]$ws.Range("D3").Select()
